$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.848.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.724.27'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.36'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -5.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4898'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +5.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3514'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.72'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07235'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.050'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.15%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.84'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.867'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.717.99'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.799'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '86.50'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -6.26%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.50'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.713'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '26.895.09'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.94'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.050'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.43'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.87'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.914.80'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.064'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.78'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.044'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09302'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.578'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.353'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.05876'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02173'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.423'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.50%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.743'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.65%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1980'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5968'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.114'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.405'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.75'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.575'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5602'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '119.57'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.833'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06629'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.093'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.96%  '

Write-Output "Applied 98 cell updates"
